$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) of the template data row (row 2) down into the new rows 63-70
$ws.Range("A2:V2").Copy()
$ws.Range("A63:V70").PasteSpecial(-4122)

# Row 63
$ws.Cells.Item(63,1).Value = 62
$ws.Cells.Item(63,2).Value = "thailand"
$ws.Cells.Item(63,3).Value = "thai-league-1"
$ws.Cells.Item(63,4).Value = "2023-2024"
$ws.Cells.Item(63,5).Value = 45226.58333333334
$ws.Cells.Item(63,6).Value = "Police Tero"
$ws.Cells.Item(63,7).Value = 1
$ws.Cells.Item(63,8).Value = "Khonkaen Utd."
$ws.Cells.Item(63,9).Value = 3
$ws.Cells.Item(63,10).Value = 1.71
$ws.Cells.Item(63,11).Value = "26/10/2023 01:42"
$ws.Cells.Item(63,12).Value = 1.7
$ws.Cells.Item(63,13).Value = "27/10/2023 13:51"
$ws.Cells.Item(63,14).Value = 3.97
$ws.Cells.Item(63,15).Value = "26/10/2023 01:42"
$ws.Cells.Item(63,16).Value = 4.21
$ws.Cells.Item(63,17).Value = "27/10/2023 13:53"
$ws.Cells.Item(63,18).Value = 4.51
$ws.Cells.Item(63,19).Value = "26/10/2023 01:42"
$ws.Cells.Item(63,20).Value = 4.5
$ws.Cells.Item(63,21).Value = "27/10/2023 13:51"
$ws.Cells.Item(63,22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/police-tero-khonkaen-united/hMX6ywA1/"

# Row 64
$ws.Cells.Item(64,1).Value = 63
$ws.Cells.Item(64,2).Value = "thailand"
$ws.Cells.Item(64,3).Value = "thai-league-1"
$ws.Cells.Item(64,4).Value = "2023-2024"
$ws.Cells.Item(64,5).Value = 45226.625
$ws.Cells.Item(64,6).Value = "Sukhothai"
$ws.Cells.Item(64,7).Value = 2
$ws.Cells.Item(64,8).Value = "Port MTI FC"
$ws.Cells.Item(64,9).Value = 1
$ws.Cells.Item(64,10).Value = 5.4
$ws.Cells.Item(64,11).Value = "26/10/2023 07:42"
$ws.Cells.Item(64,12).Value = 5.08
$ws.Cells.Item(64,13).Value = "27/10/2023 14:59"
$ws.Cells.Item(64,14).Value = 4.36
$ws.Cells.Item(64,15).Value = "26/10/2023 07:42"
$ws.Cells.Item(64,16).Value = 5.08
$ws.Cells.Item(64,17).Value = "27/10/2023 14:59"
$ws.Cells.Item(64,18).Value = 1.5
$ws.Cells.Item(64,19).Value = "26/10/2023 07:42"
$ws.Cells.Item(64,20).Value = 1.52
$ws.Cells.Item(64,21).Value = "27/10/2023 14:59"
$ws.Cells.Item(64,22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/sukhothai-port-mti-fc/bJhSvCGg/"

# Row 65
$ws.Cells.Item(65,1).Value = 64
$ws.Cells.Item(65,2).Value = "thailand"
$ws.Cells.Item(65,3).Value = "thai-league-1"
$ws.Cells.Item(65,4).Value = "2023-2024"
$ws.Cells.Item(65,5).Value = 45227.54166666666
$ws.Cells.Item(65,6).Value = "Chiangrai Utd"
$ws.Cells.Item(65,7).Value = 1
$ws.Cells.Item(65,8).Value = "Chonburi"
$ws.Cells.Item(65,9).Value = 1
$ws.Cells.Item(65,10).Value = 2.11
$ws.Cells.Item(65,11).Value = "27/10/2023 14:13"
$ws.Cells.Item(65,12).Value = 2.64
$ws.Cells.Item(65,13).Value = "28/10/2023 12:53"
$ws.Cells.Item(65,14).Value = 3.8
$ws.Cells.Item(65,15).Value = "27/10/2023 14:13"
$ws.Cells.Item(65,16).Value = 3.87
$ws.Cells.Item(65,17).Value = "28/10/2023 12:53"
$ws.Cells.Item(65,18).Value = 2.98
$ws.Cells.Item(65,19).Value = "27/10/2023 14:13"
$ws.Cells.Item(65,20).Value = 2.46
$ws.Cells.Item(65,21).Value = "28/10/2023 12:53"
$ws.Cells.Item(65,22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/chiangrai-utd-chonburi/Gf3Nuh1m/"

# Row 66
$ws.Cells.Item(66,1).Value = 65
$ws.Cells.Item(66,2).Value = "thailand"
$ws.Cells.Item(66,3).Value = "thai-league-1"
$ws.Cells.Item(66,4).Value = "2023-2024"
$ws.Cells.Item(66,5).Value = 45227.58333333334
$ws.Cells.Item(66,6).Value = "Lamphun Warrior"
$ws.Cells.Item(66,7).Value = 2
$ws.Cells.Item(66,8).Value = "Nakhon Pathom"
$ws.Cells.Item(66,9).Value = 0
$ws.Cells.Item(66,10).Value = 2.28
$ws.Cells.Item(66,11).Value = "27/10/2023 14:13"
$ws.Cells.Item(66,12).Value = 1.76
$ws.Cells.Item(66,13).Value = "28/10/2023 13:22"
$ws.Cells.Item(66,14).Value = 3.73
$ws.Cells.Item(66,15).Value = "27/10/2023 14:13"
$ws.Cells.Item(66,16).Value = 4.03
$ws.Cells.Item(66,17).Value = "28/10/2023 13:25"
$ws.Cells.Item(66,18).Value = 2.74
$ws.Cells.Item(66,19).Value = "27/10/2023 14:13"
$ws.Cells.Item(66,20).Value = 4.3
$ws.Cells.Item(66,21).Value = "28/10/2023 13:25"
$ws.Cells.Item(66,22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/lamphun-warrior-nakhon-pathom/W4xnat28/"

# Row 67
$ws.Cells.Item(67,1).Value = 66
$ws.Cells.Item(67,2).Value = "thailand"
$ws.Cells.Item(67,3).Value = "thai-league-1"
$ws.Cells.Item(67,4).Value = "2023-2024"
$ws.Cells.Item(67,5).Value = 45227.625
$ws.Cells.Item(67,6).Value = "Uthai Thani"
$ws.Cells.Item(67,7).Value = 4
$ws.Cells.Item(67,8).Value = "Muang Thong Utd"
$ws.Cells.Item(67,9).Value = 2
$ws.Cells.Item(67,10).Value = 3.63
$ws.Cells.Item(67,11).Value = "27/10/2023 14:13"
$ws.Cells.Item(67,12).Value = 2.87
$ws.Cells.Item(67,13).Value = "28/10/2023 14:57"
$ws.Cells.Item(67,14).Value = 3.79
$ws.Cells.Item(67,15).Value = "27/10/2023 14:13"
$ws.Cells.Item(67,16).Value = 3.76
$ws.Cells.Item(67,17).Value = "28/10/2023 14:56"
$ws.Cells.Item(67,18).Value = 1.93
$ws.Cells.Item(67,19).Value = "27/10/2023 14:13"
$ws.Cells.Item(67,20).Value = 2.33
$ws.Cells.Item(67,21).Value = "28/10/2023 14:57"
$ws.Cells.Item(67,22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/uthai-thani-muang-thong-utd/baY2xJfe/"

# Row 68
$ws.Cells.Item(68,1).Value = 67
$ws.Cells.Item(68,2).Value = "thailand"
$ws.Cells.Item(68,3).Value = "thai-league-1"
$ws.Cells.Item(68,4).Value = "2023-2024"
$ws.Cells.Item(68,5).Value = 45228.47916666666
$ws.Cells.Item(68,6).Value = "Trat FC"
$ws.Cells.Item(68,7).Value = 1
$ws.Cells.Item(68,8).Value = "Prachuap"
$ws.Cells.Item(68,9).Value = 0
$ws.Cells.Item(68,10).Value = 2.17
$ws.Cells.Item(68,11).Value = "28/10/2023 07:43"
$ws.Cells.Item(68,12).Value = 2.73
$ws.Cells.Item(68,13).Value = "29/10/2023 11:07"
$ws.Cells.Item(68,14).Value = 3.69
$ws.Cells.Item(68,15).Value = "28/10/2023 07:43"
$ws.Cells.Item(68,16).Value = 3.74
$ws.Cells.Item(68,17).Value = "29/10/2023 11:27"
$ws.Cells.Item(68,18).Value = 2.94
$ws.Cells.Item(68,19).Value = "28/10/2023 07:43"
$ws.Cells.Item(68,20).Value = 2.41
$ws.Cells.Item(68,21).Value = "29/10/2023 11:07"
$ws.Cells.Item(68,22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/trat-fc-prachuap/Aeyjb0HE/"

# Row 69
$ws.Cells.Item(69,1).Value = 68
$ws.Cells.Item(69,2).Value = "thailand"
$ws.Cells.Item(69,3).Value = "thai-league-1"
$ws.Cells.Item(69,4).Value = "2023-2024"
$ws.Cells.Item(69,5).Value = 45228.5
$ws.Cells.Item(69,6).Value = "Buriram"
$ws.Cells.Item(69,7).Value = 0
$ws.Cells.Item(69,8).Value = "Pathum United"
$ws.Cells.Item(69,9).Value = 0
$ws.Cells.Item(69,10).Value = 1.51
$ws.Cells.Item(69,11).Value = "27/10/2023 21:43"
$ws.Cells.Item(69,12).Value = 1.68
$ws.Cells.Item(69,13).Value = "29/10/2023 11:55"
$ws.Cells.Item(69,14).Value = 4.3
$ws.Cells.Item(69,15).Value = "27/10/2023 21:43"
$ws.Cells.Item(69,16).Value = 4.21
$ws.Cells.Item(69,17).Value = "29/10/2023 11:55"
$ws.Cells.Item(69,18).Value = 5.35
$ws.Cells.Item(69,19).Value = "27/10/2023 21:43"
$ws.Cells.Item(69,20).Value = 4.63
$ws.Cells.Item(69,21).Value = "29/10/2023 11:55"
$ws.Cells.Item(69,22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/buriram-united-f-c-pathum-united/2iaJtYos/"

# Row 70
$ws.Cells.Item(70,1).Value = 69
$ws.Cells.Item(70,2).Value = "thailand"
$ws.Cells.Item(70,3).Value = "thai-league-1"
$ws.Cells.Item(70,4).Value = "2023-2024"
$ws.Cells.Item(70,5).Value = 45228.54166666666
$ws.Cells.Item(70,6).Value = "Bangkok Utd"
$ws.Cells.Item(70,7).Value = 3
$ws.Cells.Item(70,8).Value = "Ratchaburi"
$ws.Cells.Item(70,9).Value = 1
$ws.Cells.Item(70,10).Value = 1.51
$ws.Cells.Item(70,11).Value = "27/10/2023 21:43"
$ws.Cells.Item(70,12).Value = 1.46
$ws.Cells.Item(70,13).Value = "29/10/2023 12:51"
$ws.Cells.Item(70,14).Value = 4.3
$ws.Cells.Item(70,15).Value = "27/10/2023 21:43"
$ws.Cells.Item(70,16).Value = 4.39
$ws.Cells.Item(70,17).Value = "29/10/2023 12:56"
$ws.Cells.Item(70,18).Value = 6.12
$ws.Cells.Item(70,19).Value = "27/10/2023 21:43"
$ws.Cells.Item(70,20).Value = 7.19
$ws.Cells.Item(70,21).Value = "29/10/2023 12:56"
$ws.Cells.Item(70,22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/bangkok-utd-ratchaburi/WCWAzcP7/"
